# Applies the "Extraccion total turbobeads" post-meeting adjustments.
#
# Each edit is scoped to the specific paragraph it belongs to (via the
# Paragraphs collection) so that identical phrases appearing elsewhere in
# the document (e.g. multiple "Espera de 10 minutos." or
# "... del deepwell 10 veces." steps) are not accidentally touched.

$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $findText, $replaceText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, `
                           $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "WARNING: paragraph $paraIndex find failed for [$findText]"
    }
    return $ok
}

# --- Paragraph 67: ELUTION_FINAL_VOLUME_PER_SAMPLE bullet -------------------
# Runs " Volumen en " / "uL de elution" / " que sera transferido a la placa
# final." are authored as one contiguous run; re-assert the text in place so
# the fragments collapse into a single run, matching the cleaned-up markup.
Replace-InParagraph 67 " Volumen en μL de elution que será transferido a la placa final." `
                        " Volumen en μL de elution que será transferido a la placa final."

# --- Paragraph 76: "Se mueven 300 uL (x8) ..." bullet (PASO 1) -------------
Replace-InParagraph 76 "Se mueven 300 μL (x8) del canal correspondiente del reservorio multicanal a cada una de las muestras." `
                        "Se mueven 300 μL (x8) del canal correspondiente del reservorio multicanal a cada una de las muestras."

# --- Paragraph 96: "Se mueven 180 uL (x8), tantas veces..." (PASO 6) -------
Replace-InParagraph 96 " (x8), tantas veces como sea necesario para remover todo el sobrenadante," `
                        " (x8), tantas veces como sea necesario para remover todo el sobrenadante,"

# --- Paragraph 104: "Se resuspenden 180 uL del deepwell 10 veces." (PASO 8) -
Replace-InParagraph 104 " del deepwell 10 veces." " del deepwell 20 veces."

# --- Paragraph 120: same bullet repeated under PASO 12 ----------------------
Replace-InParagraph 120 " del deepwell 10 veces." " del deepwell 20 veces."

# --- Paragraph 131: "Espera de 10 minutos." under PASO 15 (Allow dry) ------
Replace-InParagraph 131 "Espera de 10 minutos." "Espera de 15 minutos."

# --- Paragraph 137: "Se mueven 75 uL ..." under PASO 17 (Transfer elution) --
Replace-InParagraph 137 "Se mueven 75 " "Se mueven 50 "

# --- Paragraph 138: "Se resuspenden 75 uL del deepwell 5 veces." -----------
Replace-InParagraph 138 "Se resuspenden 75 " "Se resuspenden 50 "
Replace-InParagraph 138 " del deepwell 5 veces." " del deepwell 20 veces."

# --- Paragraph 141: "Espera de 5 minuto." under PASO 18 ---------------------
Replace-InParagraph 141 " minuto." " minutos."

# --- Paragraph 144: "Espera de 10 minutos." under PASO 19 -------------------
Replace-InParagraph 144 "Espera de 10 minutos." "Espera de 5 minutos."

# --- Paragraph 149: "Se mueven 50 uL ..." under PASO 20 ---------------------
Replace-InParagraph 149 "Se mueven 50 " "Se mueven 45 "

Write-Output "done"
